$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.448.10'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '2.101.68'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5214'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4538'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.69'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +16.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08897'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.179'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.05'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '2.102.99'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.804'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.990'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001141'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.006'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06648'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.302'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").Value = '30.541.78'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.341'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").Value = '2.359.89'
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.522'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.201'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.648'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.415'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.950'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.782'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02574'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2298'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6850'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.246'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.317'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6346'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.663'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("E48").Value = '  +20.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.250'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.204'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '83.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
